$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.049.08"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.305.94"
$ws.Range("E3").Value = "  -0.83%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.64"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.28"
$ws.Range("E6").Value = "  -1.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.519"
$ws.Range("E7").Value = "  +2.26%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.512"
$ws.Range("E9").Value = "  -1.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.38"
$ws.Range("E10").Value = "  -0.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.92"
$ws.Range("E12").Value = "  +0.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.82"
$ws.Range("E14").Value = "  -2.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.663.48"
$ws.Range("E15").Value = "  -0.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.311.71"
$ws.Range("E16").Value = "  -2.81%  "

$ws.Range("E17").Value = "  -2.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.006.39"
$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  +0.68%  "

$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("E21").Value = "  -1.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.42"
$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.98"
$ws.Range("E23").Value = "  +1.68%  "

$ws.Range("E24").Value = "  -1.39%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("E26").Value = "  -1.27%  "

$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.21"
$ws.Range("E28").Value = "  -0.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.48"
$ws.Range("E29").Value = "  -1.08%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.10"
$ws.Range("E30").Value = "  -0.89%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.03"
$ws.Range("E31").Value = "  -0.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.34"
$ws.Range("E32").Value = "  -3.61%  "

$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("E34").Value = "  -3.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.73"
$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.80"
$ws.Range("E36").Value = "  +1.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.41"
$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0691"
$ws.Range("E38").Value = "  -0.64%  "

$ws.Range("E39").Value = "  -1.41%  "

$ws.Range("E40").Value = "  -1.74%  "

$ws.Range("E41").Value = "  +0.54%  "

$ws.Range("E42").Value = "  +0.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.998.67"
$ws.Range("E43").Value = "  -0.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0286"
$ws.Range("E44").Value = "  -1.91%  "

$ws.Range("E45").Value = "  -3.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.25"
$ws.Range("E46").Value = "  +1.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.41"
$ws.Range("E47").Value = "  -2.72%  "

$ws.Range("E48").Value = "  -3.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.76"
$ws.Range("E49").Value = "  -2.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.529.59"
$ws.Range("E50").Value = "  -0.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.81"
$ws.Range("E51").Value = "  -0.42%  "
